# Cancer Crosswalk.xlsx - add HUQ051 questionnaire crosswalk sheet
$wb = $excel.ActiveWorkbook

# --- Update the selection on PUQ110 (now just a header-row selection) ---
$puq = $wb.Worksheets.Item("PUQ110")
[void]$puq.Range("A1:B1").Select()

# --- Insert the new HUQ051 sheet right after HUQ071 ---
$after = $wb.Worksheets.Item("HUQ071")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "HUQ051"

$ws.Range("A1").Value = "Answer"
$ws.Range("B1").Value = "Coefficient"

$ws.Range("A2").Value = "None"
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "2 to 3"
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = "4 to 5"
$ws.Range("B5").Value = 3

$ws.Range("A6").Value = "6 to 7"
$ws.Range("B6").Value = 4

$ws.Range("A7").Value = "8 to 9"
$ws.Range("B7").Value = 5

$ws.Range("A8").Value = "10 to 12"
$ws.Range("B8").Value = 6

$ws.Range("A9").Value = "13 to 15"
$ws.Range("B9").Value = 7

$ws.Range("A10").Value = "16 or more"
$ws.Range("B10").Value = 8

$ws.Range("A11").Value = "Refused"
$ws.Range("B11").Value = 77

$ws.Range("A12").Value = "Don't Know"
$ws.Range("B12").Value = 99

$ws.Range("A13").Value = "Missing"
$ws.Range("B13").Value = "."

# Leave the cursor where Excel would land after typing down column B,
# and make this newly-added sheet the active tab.
[void]$ws.Range("B14").Select()
$ws.Activate()
